$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 29 ---
$ws.Range("B29").Value = "Wages Expense"
$ws.Range("D29").Formula = "=45000+195000"

# --- Row 30 ---
$ws.Range("B30").Value = "BELI kresek"
$ws.Range("D30").Value = 97000

# --- Row 31 ---
$ws.Range("B31").Value = "TRANSFER BCA"
$ws.Range("D31").Formula = "=464000+3680000+170000+230000+1082000"

# --- Row 32 ---
$ws.Range("B32").Value = "TELPON - 5224823"
$ws.Range("D32").Value = 681500

# --- Row 33 ---
$ws.Range("B33").Value = "PLN - Astar 165"
$ws.Range("D33").Formula = "=791500"

# --- Row 34 ---
$ws.Range("B34").Value = "CHEQUE RECEIVED"
$ws.Range("D34").Formula = "=1670000"

# --- Row 35 ---
$ws.Range("B35").Value = "PRIVE - andreas"
$ws.Range("D35").Value = 5000000

# --- Row 36 ---
$ws.Range("B36").Value = "A/R"
$ws.Range("C36").Formula = "=22846000"

# --- Row 37 ---
$ws.Range("B37").Value = "SALES - cash/retail"
$ws.Range("C37").Formula = "=12734975+15933525-22846000"

# --- Row 38 ---
$ws.Range("B38").Value = "SETOR KE BANK"
$ws.Range("D38").Value = 14000000

# --- Row 39 (new date) ---
$ws.Range("A39").Value = 44274
$ws.Range("B39").Value = "Wages Expense"
$ws.Range("D39").Formula = "=45000+195000"

# --- Row 40 ---
$ws.Range("B40").Value = "TRANSFER BCA"
$ws.Range("D40").Formula = "=500000+788000+415000+18600000+14500000"

# --- Row 41 ---
$ws.Range("B41").Value = "A/P"
$ws.Range("D41").Formula = "=2150000"

# --- Row 42 ---
$ws.Range("B42").Value = "A/R"
$ws.Range("C42").Formula = "=18268000+13000000+7236000+10013000"

# --- Row 43 ---
$ws.Range("B43").Value = "SALES - cash/retail"
$ws.Range("C43").Formula = "=24668525-2162525-10013000"

# --- Row 44 ---
$ws.Range("B44").Value = "SERVICE rumah MH27"
$ws.Range("D44").Value = 500000

# --- Row 45 ---
$ws.Range("B45").Value = "SELISIH - lebih"
$ws.Range("C45").Value = 500

# --- Row 46 ---
$ws.Range("B46").Value = "SETOR KE BANK"
$ws.Range("D46").Value = 23000000

# --- Row 47 (new date) ---
$ws.Range("A47").Value = 44275
$ws.Range("B47").Value = "Wages Expense"

# --- View state: frozen pane scroll position + active selection ---
$ws.Application.ActiveWindow.ScrollRow = 46
$ws.Range("C67").Select()
